$d = $word.ActiveDocument

# Locate the two anchor paragraphs by text content to be resilient to any
# paragraph-index drift. The replaced region spans from the paragraph
# containing "路由选择、维护路由表" through to (and including) the final
# paragraph containing "应用层：各种应用程序，协议" (the old simple layer-7
# line, which is always the very last occurrence of that phrase before the
# edit).
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*路由选择*维护路由表*") {
        $startPara = $p
    }
    if ($t -like "*应用层：各种应用程序*") {
        $endPara = $p
    }
}

if ($startPara -eq $null -or $endPara -eq $null) {
    Write-Host "ERROR: could not locate anchor paragraphs"
} else {
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()

    $insertPoint = $d.Range($startPara.Range.Start, $startPara.Range.Start)

    $xmlFragment = @'
<w:p><w:pPr><w:widowControl w:val="0"/><w:tabs><w:tab w:val="left" w:pos="400"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>3）路由选择、维护路由表</w:t></w:r></w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Helvetica"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Helvetica" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">4. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>传输层：提供可靠或不可靠的端到端传输</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>1)负责将来自上层应用程序的数据进行分段和重组，并将它们组合为同样地数据流形式</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>2)提供端到端的数据传输服务</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>3)工作在传输层的协议：TCP、UDP</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Helvetica"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Helvetica" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">5. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>会话层：为通信双方指定通信方式，并创建，注销会话</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>1）负责建立、管理和终止表示层实体之间的会话连接</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>2）在设备或节点之间提供会话控制</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>3）在系统之间协调通信过程</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Helvetica"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Helvetica" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">6. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>表示层：数据和信息的转换，数据压缩解压缩，加密解密</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>1）数据的解码和编码</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>2）</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>数据的</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>加密和解密</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>3）</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>数据的</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>压缩和解压</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:widowControl w:val="0"/>
        <w:tabs>
          <w:tab w:val="left" w:pos="400"/>
        </w:tabs>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>4）常见的标准：ASCII、JPEG</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Helvetica"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Helvetica" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Helvetica" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">7. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="PingFang SC" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>应用层：各种应用程序，协议</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Helvetica" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Helvetica" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Helvetica" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>1）为应用软件提供接口，使应用程序能够使用网络服务</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Helvetica" w:hint="eastAsia"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>2）常见应用层协议：http、ftp、dns、telnet</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p><w:p><w:r><w:t>ZZZDUMMYZZZ</w:t></w:r></w:p>
'@

    $insertPoint.InsertXML($xmlFragment)

    # Remove the dummy run of text that absorbed the trailing-paragraph-mark
    # merge behaviour of InsertXML.
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $found = $find.Execute("ZZZDUMMYZZZ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
    Write-Host "Dummy found/replaced:" $found

    Write-Host "Final paragraph count:" $d.Paragraphs.Count
}
